$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E6").Value = "Problem 17"
$ws.Range("F7").Value = "string initlialization"

$ws.Range("F7").Select()
